# The three images in this document (two copies of the Pearson logo in
# the footers, and the BTEC logo in the header) had their internal
# drawing names swapped:
#   - Pearson logo (footer1.xml, docPr id="3")  image1.png -> image2.png
#   - Pearson logo (footer2.xml, docPr id="2")  image1.png -> image2.png
#   - BTEC logo    (header1.xml, docPr id="1")  image2.jpg -> image1.jpg
#
# Word's InlineShape object model does not expose a writable "Name"
# property that reaches both the <wp:docPr> AND the nested
# <pic:cNvPr> elements, so we round-trip the package through
# Document.WordOpenXML (the flat-OPC single-XML view of the whole
# document, headers/footers included) and patch the `name="..."`
# attributes there, then write the edited XML back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# --- Footer 1 (id="3"): Pearson logo docPr -----------------------------
$xml = $xml.Replace(
    'docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"',
    'docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"')

# --- Footer 2 (id="2"): Pearson logo docPr -----------------------------
$xml = $xml.Replace(
    'docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"',
    'docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"')

# --- Both Pearson logo <pic:cNvPr> (id="0") -----------------------------
# Both footers use the identical cNvPr text, so one global Replace
# rewrites both occurrences (matching the diff, which touches both).
$xml = $xml.Replace(
    'pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"',
    'pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"')

# --- Header 1: BTEC logo docPr and cNvPr --------------------------------
$xml = $xml.Replace(
    'docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"',
    'docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"')

$xml = $xml.Replace(
    'pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"',
    'pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"')

$d.WordOpenXML = $xml
